$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value of 45190 (2023-09-21) for every
# data row (rows 2-261). Update it to 45192 (2023-09-23) for all of them.
$ws.Range("C2:C261").Value = 45192
